$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# Page margins / footer distance tweaks to match target sectPr
$section.PageSetup.TopMargin = 36.2
$section.PageSetup.FooterDistance = 14.6

# Turn on line numbering for the section (lnNumType countBy=1)
$section.PageSetup.LineNumbering.Active = $true
$section.PageSetup.LineNumbering.CountBy = 1

# Add page number field to the footer, right aligned
$footer = $section.Footers.Item(1)
$footer.PageNumbers.Add(2)
$footer.Range.ParagraphFormat.Alignment = 2
